# Changed Test Data for LV Activities - 16 Sep 2024
#
# - Contact sheet: replace the "Test External" / "StandardTestCompany"
#   test row with "Activity Test External Contact" / "ActivityCompany".
# - Move the active/selected tab from "MoreAttendees" to "Contact",
#   and update the Contact sheet's selection to A2:B2.

$wb = $excel.ActiveWorkbook

$contact = $wb.Worksheets.Item("Contact")
$contact.Range("A2").Value = "Activity Test External Contact"
$contact.Range("B2").Value = "ActivityCompany"

$contact.Activate()
$contact.Range("A2:B2").Select()
